$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold plain-text numbers (thousands "." separators, no
# true numeric typing). Any single-dot-decimal-looking string assigned via
# .Value would otherwise be auto-coerced by Excel into a real number and
# lose its exact text (trailing zeros, float rounding). Force those specific
# cells to Text format first so the literal string is preserved.

$ws.Range("D2").Value = '28.592.32'
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").Value = '1.913.76'
$ws.Range("E3").Value = '  +5.70%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.83'
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5055'
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3957'
$ws.Range("E8").Value = '  +1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09815'
$ws.Range("E9").Value = '  +1.10%  '
$ws.Range("E10").Value = '  +5.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.37'
$ws.Range("E11").Value = '  +3.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.548'
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("E13").Value = '  +3.82%  '
$ws.Range("D14").Value = '1.914.22'
$ws.Range("E14").Value = '  +5.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.580'
$ws.Range("E15").Value = '  +4.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001142'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.14'
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06659'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.14'
$ws.Range("E20").Value = '  +5.83%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.313'
$ws.Range("E22").Value = '  +6.86%  '
$ws.Range("D23").Value = '28.652.17'
$ws.Range("E23").Value = '  +2.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.46'
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.281'
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.746'
$ws.Range("E26").Value = '  +15.48%  '
$ws.Range("D27").Value = '2.135.91'
$ws.Range("E27").Value = '  +5.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.31'
$ws.Range("E28").Value = '  +3.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.33'
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.95'
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.105'
$ws.Range("E31").Value = '  +6.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1074'
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.753'
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.646'
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.860'
$ws.Range("E35").Value = '  +11.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06806'
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("E37").Value = '  +5.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.275'
$ws.Range("E38").Value = '  +9.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2233'
$ws.Range("E39").Value = '  +4.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.78'
$ws.Range("E40").Value = '  +4.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.099'
$ws.Range("E42").Value = '  +4.69%  '
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.73'
$ws.Range("E45").Value = '  +4.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6093'
$ws.Range("E46").Value = '  +3.66%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.827'
$ws.Range("E47").Value = '  +3.83%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.284'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.045'
$ws.Range("E49").Value = '  +5.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.19'
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("E51").Value = '  +3.22%  '
